$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New data for rows 2-12 (columns A-I)
$data = @(
  @(0,  "5008633-12.2022.8.21.4001", "5002299-64.2019.8.21.4001", "Originário",  "Sem dados de processo originário 2", "Nulo", "Sem dados de processo originário 3", "Nulo", "14/10/2022"),
  @(1,  "5006127-55.2022.8.21.0059", "5001247-88.2020.8.21.0059", "Originário",  "Sem dados de processo originário 2", "Nulo", "Sem dados de processo originário 3", "Nulo", "19/09/2022"),
  @(2,  "5157435-21.2022.8.21.0001", "5032648-85.2020.8.21.0001", "Originário",  "Sem dados de processo originário 2", "Nulo", "Sem dados de processo originário 3", "Nulo", "05/09/2022"),
  @(3,  "5026797-60.2023.8.21.0001", "5030835-67.2013.8.21.0001", "Originário",  "Sem dados de processo originário 2", "Nulo", "Sem dados de processo originário 3", "Nulo", "16/02/2023"),
  @(4,  "5035699-12.2017.8.21.0001", "0161621-51.2017.8.21.0001", "Digitalizado","Sem dados de processo originário 2", "Nulo", "Sem dados de processo originário 3", "Nulo", "31/08/2017"),
  @(5,  "5000273-46.2011.8.21.0001", "0420275-57.2011.8.21.0001", "Digitalizado","Sem dados de processo originário 2", "Nulo", "Sem dados de processo originário 3", "Nulo", "07/12/2011"),
  @(6,  "5000274-31.2011.8.21.0001", "0413893-48.2011.8.21.0001", "Digitalizado","Sem dados de processo originário 2", "Nulo", "Sem dados de processo originário 3", "Nulo", "01/12/2011"),
  @(7,  "5004401-12.2011.8.21.0001", "0376983-22.2011.8.21.0001", "Digitalizado","Sem dados de processo originário 2", "Nulo", "Sem dados de processo originário 3", "Nulo", "07/11/2011"),
  @(8,  "5000256-73.2012.8.21.0001", "0023506-26.2012.8.21.0001", "Digitalizado","Sem dados de processo originário 2", "Nulo", "Sem dados de processo originário 3", "Nulo", "24/01/2012"),
  @(9,  "5000271-42.2012.8.21.0001", "0351565-48.2012.8.21.0001", "Digitalizado","Sem dados de processo originário 2", "Nulo", "Sem dados de processo originário 3", "Nulo", "31/10/2012"),
  @(10, "5000293-03.2012.8.21.0001", "0342323-65.2012.8.21.0001", "Digitalizado","Sem dados de processo originário 2", "Nulo", "Sem dados de processo originário 3", "Nulo", "24/10/2012")
)

$row = 2
foreach ($r in $data) {
  $ws.Cells.Item($row, 1).Value = $r[0]
  $ws.Cells.Item($row, 2).Value = $r[1]
  $ws.Cells.Item($row, 3).Value = $r[2]
  $ws.Cells.Item($row, 4).Value = $r[3]
  $ws.Cells.Item($row, 5).Value = $r[4]
  $ws.Cells.Item($row, 6).Value = $r[5]
  $ws.Cells.Item($row, 7).Value = $r[6]
  $ws.Cells.Item($row, 8).Value = $r[7]

  # The distribution-date column is plain text in the source data. Excel's
  # automatic type inference turns some DD/MM/YYYY strings into real dates
  # whenever the day portion is <= 12 (so it also parses as a month). Force
  # those specific cells to Text format first so the value is kept verbatim.
  $dateCell = $ws.Cells.Item($row, 9)
  $day = [int]($r[8].Split('/')[0])
  if ($day -le 12) {
    $dateCell.NumberFormat = "@"
  }
  $dateCell.Value = $r[8]

  $row++
}

# Remove the now-obsolete rows 13-19 entirely
$ws.Range("A13:I19").EntireRow.Delete()
